$wb = $excel.ActiveWorkbook

# Update the Date value on the Metadata sheet
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-07-25T07:22:51+00:00"

# Update the System URI on the Include #0 sheet
$wsInc0 = $wb.Worksheets.Item("Include #0")
$wsInc0.Range("B4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R45-FonctionQualifiee/FHIR/TRE-R45-FonctionQualifiee"

# Update the System URI on the Include #1 sheet
$wsInc1 = $wb.Worksheets.Item("Include #1")
$wsInc1.Range("B4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R01-EnsembleSavoirFaire-CISIS/FHIR/TRE-R01-EnsembleSavoirFaire-CISIS"
